$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.458.29'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Value = '2.068.93'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '232.21'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').Value = '0.631'
$ws.Range('E6').Value = '  +1.35%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '57.34'
$ws.Range('E8').Value = '  -2.20%  '
$ws.Range('D9').Value = '0.389'
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('D10').Value = '0.0777'
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('E11').Value = '  +1.54%  '
$ws.Range('D12').Value = '14.82'
$ws.Range('E12').Value = '  +0.56%  '
$ws.Range('D13').Value = '2.373.57'
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('D14').Value = '20.85'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').Value = '0.761'
$ws.Range('E15').Value = '  -1.64%  '
$ws.Range('D16').Value = '5.31'
$ws.Range('E16').Value = '  -0.98%  '
$ws.Range('D17').Value = '2.067.12'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '37.358.11'
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('D19').Value = '70.41'
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('E20').Value = '  -2.58%  '
$ws.Range('D21').Value = '0.0₃0826'
$ws.Range('E21').Value = '  -0.97%  '
$ws.Range('D22').Value = '227.89'
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('D26').Value = '9.64'
$ws.Range('E26').Value = '  +6.49%  '
$ws.Range('D27').Value = '169.63'
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('E28').Value = '  -3.69%  '
$ws.Range('D29').Value = '19.45'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = '1.38'
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('D31').Value = '0.123'
$ws.Range('E31').Value = '  +1.42%  '
$ws.Range('D32').Value = '4.61'
$ws.Range('E32').Value = '  -1.68%  '
$ws.Range('E33').Value = '  -0.39%  '
$ws.Range('D34').Value = '4.61'
$ws.Range('E34').Value = '  -1.24%  '
$ws.Range('D35').Value = '2.47'
$ws.Range('E35').Value = '  -0.78%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').Value = '3.32'
$ws.Range('E37').Value = '  -1.84%  '
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').Value = '5.26'
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('D40').Value = '0.0230'
$ws.Range('E40').Value = '  +6.82%  '
$ws.Range('D41').Value = '99.58'
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('D42').Value = '2.91'
$ws.Range('E42').Value = '  +0.81%  '
$ws.Range('D43').Value = '1.20'
$ws.Range('E43').Value = '  +3.91%  '
$ws.Range('D44').Value = '0.0950'
$ws.Range('E44').Value = '  -2.34%  '
$ws.Range('D45').Value = '1.476.19'
$ws.Range('E45').Value = '  +2.50%  '
$ws.Range('D46').Value = '16.70'
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('E47').Value = '  -1.63%  '
$ws.Range('E48').Value = '  -2.05%  '
$ws.Range('D49').Value = '3.92'
$ws.Range('E49').Value = '  -6.22%  '
$ws.Range('E50').Value = '  -2.15%  '
$ws.Range('D51').Value = '2.257.17'
$ws.Range('E51').Value = '  -0.37%  '
